# Apply cryptocurrency price/volume updates as described by the commit diff.
# The workbook is already open; operate on the active sheet (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must stay plain text even if it
# looks like a number (Excel would otherwise auto-convert values
# such as "1.003" into the numeric value 1.003). We briefly mark
# the cell as Text, set the value, then clear the formatting again
# so the cell keeps no explicit style, just like the source file.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# Row 2
$ws.Range("D2").Value = "26.004.78"
$ws.Range("E2").Value = "  +0.32%  "
# Row 3
$ws.Range("D3").Value = "1.641.59"
$ws.Range("E3").Value = "  -0.01%  "
# Row 4
Set-TextValue "D4" "1.003"
$ws.Range("E4").Value = "  -0.30%  "
# Row 5
Set-TextValue "D5" "215.17"
$ws.Range("E5").Value = "  -0.10%  "
# Row 6
Set-TextValue "D6" "0.5078"
$ws.Range("E6").Value = "  -0.25%  "
# Row 7
Set-TextValue "D7" "1.002"
$ws.Range("E7").Value = "  -0.35%  "
# Row 9
Set-TextValue "D9" "0.06364"
$ws.Range("E9").Value = "  -0.48%  "
# Row 10
Set-TextValue "D10" "19.87"
$ws.Range("E10").Value = "  +1.60%  "
# Row 11
Set-TextValue "D11" "0.07734"
$ws.Range("E11").Value = "  -0.47%  "
# Row 12
Set-TextValue "D12" "4.302"
$ws.Range("E12").Value = "  -0.06%  "
# Row 13
$ws.Range("D13").Value = "1.636.35"
$ws.Range("E13").Value = "  -0.50%  "
# Row 14
Set-TextValue "D14" "0.5478"
$ws.Range("E14").Value = "  +0.40%  "
# Row 15
$ws.Range("D15").Value = "0.0₅7755"
$ws.Range("E15").Value = "  -1.25%  "
# Row 16
Set-TextValue "D16" "64.42"
$ws.Range("E16").Value = "  -0.39%  "
# Row 17
$ws.Range("D17").Value = "26.032.05"
$ws.Range("E17").Value = "  +0.18%  "
# Row 18
Set-TextValue "D18" "1.003"
$ws.Range("E18").Value = "  -0.29%  "
# Row 19
Set-TextValue "D19" "197.47"
$ws.Range("E19").Value = "  -0.18%  "
# Row 20
Set-TextValue "D20" "4.466"
$ws.Range("E20").Value = "  +0.64%  "
# Row 21
Set-TextValue "D21" "9.980"
$ws.Range("E21").Value = "  +0.07%  "
# Row 22
Set-TextValue "D22" "6.138"
$ws.Range("E22").Value = "  +1.71%  "
# Row 23
Set-TextValue "D23" "1.003"
$ws.Range("E23").Value = "  -0.46%  "
# Row 24
Set-TextValue "D24" "1.896"
$ws.Range("E24").Value = "  +1.01%  "
# Row 25
Set-TextValue "D25" "142.67"
$ws.Range("E25").Value = "  +1.32%  "
# Row 26
Set-TextValue "D26" "0.1258"
$ws.Range("E26").Value = "  +9.69%  "
# Row 27
Set-TextValue "D27" "6.879"
$ws.Range("E27").Value = "  -0.39%  "
# Row 28
Set-TextValue "D28" "15.62"
$ws.Range("E28").Value = "  -0.69%  "
# Row 29
Set-TextValue "D29" "1.242"
# Row 30
Set-TextValue "D30" "0.04904"
$ws.Range("E30").Value = "  -2.22%  "
# Row 31
Set-TextValue "D31" "3.287"
$ws.Range("E31").Value = "  +0.73%  "
# Row 32
Set-TextValue "D32" "3.218"
$ws.Range("E32").Value = "  +0.92%  "
# Row 33
Set-TextValue "D33" "1.561"
$ws.Range("E33").Value = "  +1.29%  "
# Row 34
Set-TextValue "D34" "2.379"
$ws.Range("E34").Value = "  +0.62%  "
# Row 35
Set-TextValue "D35" "0.9194"
$ws.Range("E35").Value = "  +2.84%  "
# Row 36
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D36" "2.567"
$ws.Range("E36").Value = "  -0.90%  "
# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D37" "0.5562"
$ws.Range("E37").Value = "  +0.86%  "
# Row 38
$ws.Range("D38").Value = "1.105.73"
$ws.Range("E38").Value = "  -2.65%  "
# Row 39
Set-TextValue "D39" "0.01569"
$ws.Range("E39").Value = "  +0.82%  "
# Row 40
Set-TextValue "D40" "1.002"
$ws.Range("E40").Value = "  -0.42%  "
# Row 41
Set-TextValue "D41" "5.616"
$ws.Range("E41").Value = "  -0.36%  "
# Row 42
Set-TextValue "D42" "0.8052"
$ws.Range("E42").Value = "  -1.51%  "
# Row 43
Set-TextValue "D43" "98.76"
$ws.Range("E43").Value = "  -1.06%  "
# Row 44
$ws.Range("D44").Value = "0.0₈122"
$ws.Range("E44").Value = "  -4.76%  "
# Row 45
$ws.Range("D45").Value = "1.781.38"
$ws.Range("E45").Value = "  +0.07%  "
# Row 46
Set-TextValue "D46" "0.4539"
$ws.Range("E46").Value = "  +0.14%  "
# Row 47
Set-TextValue "D47" "55.40"
$ws.Range("E47").Value = "  +0.84%  "
# Row 48
$ws.Range("E48").Value = "  -0.22%  "
# Row 49
Set-TextValue "D49" "0.05195"
$ws.Range("E49").Value = "  +2.12%  "
# Row 50
Set-TextValue "D50" "7.580"
$ws.Range("E50").Value = "  +2.04%  "
# Row 51
Set-TextValue "D51" "1.003"
$ws.Range("E51").Value = "  -0.26%  "
